$d = $word.ActiveDocument

# Helper: force a structural run-rebuild even when the replacement text is
# identical to the concatenated original text (the engine no-ops a Text
# assignment that doesn't change the visible text, so we bounce the text
# through a unique placeholder first, then set the real text on the 2nd pass;
# that second pass still coalesces touching same-formatted sibling runs and
# drops any wholly-enclosed <w:proofErr/> markers).
function Replace-Text($searchText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }
    if ($rng.Text -ceq $newText) {
        $placeholder = "zzPLACEHOLDERzz" + [guid]::NewGuid().ToString("N")
        $rng.Text = $placeholder
        $rng2 = $d.Content
        $rng2.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $rng2.Text = $newText
    } else {
        $rng.Text = $newText
    }
    return $true
}

# 1) "Format String: 64 bit application 3" - merge 3 runs, drop proofErr pair
Replace-Text "Format String: 64 bit application 3" "Format String: 64 bit application 3"

# 2) Insert " chạy lệnh c và" after "debug" in the gdb step instructions
Replace-Text "Sau đó bên cửa sổ debug ấn enter đến khi nào format có chứa “" "Sau đó bên cửa sổ debug chạy lệnh c và ấn enter đến khi nào format có chứa “"

# 3) "b delete" - merge 2 runs, drop proofErr pair (trailing proofErr sits at
#    the very end of the paragraph, outside addressable text, see below)
Replace-Text "b delete" "b delete"

# 4) "# function to write value3 to value4" - merge 3 runs, drop proofErr pair
Replace-Text "# function to write value3 to value4" "# function to write value3 to value4"

# 5) "# setup" - merge 2 runs, drop proofErr pair (trailing proofErr again)
Replace-Text "# setup" "# setup"

# 6) "# write_addr" - merge 3 runs, drop proofErr pair
Replace-Text "# write_addr" "# write_addr"

# 7) "# create block size 100" - merge 3 runs, drop proofErr pair
Replace-Text "# create block size 100" "# create block size 100"

# 8) "# trigger bug, leak stack, libc" - merge 3 runs, drop proofErr pair
Replace-Text "# trigger bug, leak stack, libc" "# trigger bug, leak stack, libc"

# 9) "#pause()" occurs twice - merge 3 runs, drop proofErr pair, each time
Replace-Text "#pause()" "#pause()"
Replace-Text "#pause()" "#pause()"
